$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Matriz_Resultados")
$ws2 = $wb.Worksheets.Item("P_valores")
$ws3 = $wb.Worksheets.Item("Estadisticos_DM")
$ws4 = $wb.Worksheets.Item("Resumen")

# Matriz_Resultados: 10 cell updates
$ws1.Range("J2").Value = 0
$ws1.Range("G4").Value = 0
$ws1.Range("H4").Value = 0
$ws1.Range("D7").Value = 0
$ws1.Range("J7").Value = 0
$ws1.Range("D8").Value = 0
$ws1.Range("J8").Value = 0
$ws1.Range("B10").Value = 0
$ws1.Range("G10").Value = 0
$ws1.Range("H10").Value = 0

# P_valores: 72 cell updates
$ws2.Range("C2").Value = [double]"8.408603595633224E-06"
$ws2.Range("D2").Value = [double]"1.31779794423359E-05"
$ws2.Range("E2").Value = 0.0001865587530940349
$ws2.Range("F2").Value = [double]"4.184012974528528E-06"
$ws2.Range("G2").Value = [double]"4.64420659240794E-06"
$ws2.Range("H2").Value = [double]"3.202235528076969E-06"
$ws2.Range("I2").Value = [double]"3.29857043901427E-08"
$ws2.Range("J2").Value = 0.007038606957391558
$ws2.Range("B3").Value = [double]"8.408603595633224E-06"
$ws2.Range("D3").Value = [double]"1.635874163907758E-06"
$ws2.Range("E3").Value = [double]"8.267706480546977E-06"
$ws2.Range("F3").Value = [double]"1.854678043455849E-05"
$ws2.Range("G3").Value = 0.0003993917247824275
$ws2.Range("H3").Value = 0.0005318933665054981
$ws2.Range("I3").Value = 0.01630194743413038
$ws2.Range("J3").Value = [double]"9.100368569825434E-08"
$ws2.Range("B4").Value = [double]"1.31779794423359E-05"
$ws2.Range("C4").Value = [double]"1.635874163907758E-06"
$ws2.Range("E4").Value = [double]"1.998764452570967E-05"
$ws2.Range("F4").Value = [double]"3.03569705202289E-05"
$ws2.Range("G4").Value = 0.002475749377058856
$ws2.Range("H4").Value = 0.002915019399689633
$ws2.Range("I4").Value = 0.04128972146939525
$ws2.Range("J4").Value = [double]"1.30932399367012E-07"
$ws2.Range("B5").Value = 0.0001865587530940349
$ws2.Range("C5").Value = [double]"8.267706480546977E-06"
$ws2.Range("D5").Value = [double]"1.998764452570967E-05"
$ws2.Range("F5").Value = 0.0005143821145749961
$ws2.Range("G5").Value = 0.5936081964043503
$ws2.Range("H5").Value = 0.8180907207846226
$ws2.Range("I5").Value = 0.9995570706161845
$ws2.Range("J5").Value = [double]"2.798750032018127E-05"
$ws2.Range("B6").Value = [double]"4.184012974528528E-06"
$ws2.Range("C6").Value = [double]"1.854678043455849E-05"
$ws2.Range("D6").Value = [double]"3.03569705202289E-05"
$ws2.Range("E6").Value = 0.0005143821145749961
$ws2.Range("G6").Value = [double]"1.194915978008915E-05"
$ws2.Range("H6").Value = [double]"8.844723792478248E-06"
$ws2.Range("I6").Value = [double]"1.845272026290701E-07"
$ws2.Range("J6").Value = 0.02157687694398569
$ws2.Range("B7").Value = [double]"4.64420659240794E-06"
$ws2.Range("C7").Value = 0.0003993917247824275
$ws2.Range("D7").Value = 0.002475749377058856
$ws2.Range("E7").Value = 0.5936081964043503
$ws2.Range("F7").Value = [double]"1.194915978008915E-05"
$ws2.Range("H7").Value = 0.128062116163929
$ws2.Range("I7").Value = 0.6601372073520753
$ws2.Range("J7").Value = 0.001732803998047316
$ws2.Range("B8").Value = [double]"3.202235528076969E-06"
$ws2.Range("C8").Value = 0.0005318933665054981
$ws2.Range("D8").Value = 0.002915019399689633
$ws2.Range("E8").Value = 0.8180907207846226
$ws2.Range("F8").Value = [double]"8.844723792478248E-06"
$ws2.Range("G8").Value = 0.128062116163929
$ws2.Range("I8").Value = 0.8217669694081375
$ws2.Range("J8").Value = 0.004353321525517373
$ws2.Range("B9").Value = [double]"3.29857043901427E-08"
$ws2.Range("C9").Value = 0.01630194743413038
$ws2.Range("D9").Value = 0.04128972146939525
$ws2.Range("E9").Value = 0.9995570706161845
$ws2.Range("F9").Value = [double]"1.845272026290701E-07"
$ws2.Range("G9").Value = 0.6601372073520753
$ws2.Range("H9").Value = 0.8217669694081375
$ws2.Range("J9").Value = 0.04649469443111598
$ws2.Range("B10").Value = 0.007038606957391558
$ws2.Range("C10").Value = [double]"9.100368569825434E-08"
$ws2.Range("D10").Value = [double]"1.30932399367012E-07"
$ws2.Range("E10").Value = [double]"2.798750032018127E-05"
$ws2.Range("F10").Value = 0.02157687694398569
$ws2.Range("G10").Value = 0.001732803998047316
$ws2.Range("H10").Value = 0.004353321525517373
$ws2.Range("I10").Value = 0.04649469443111598

# Estadisticos_DM: 72 cell updates
$ws3.Range("C2").Value = 6.813793025258791
$ws3.Range("D2").Value = 6.537516359248388
$ws3.Range("E2").Value = 5.022538861161706
$ws3.Range("F2").Value = 7.256117192520128
$ws3.Range("G2").Value = 7.18893327356635
$ws3.Range("H2").Value = 7.430047475841675
$ws3.Range("I2").Value = 10.86844846190304
$ws3.Range("J2").Value = 3.153776194383034
$ws3.Range("B3").Value = -6.813793025258791
$ws3.Range("D3").Value = -7.878453037160982
$ws3.Range("E3").Value = -6.824309757728702
$ws3.Range("F3").Value = -6.331545514136495
$ws3.Range("G3").Value = -4.61693374871366
$ws3.Range("H3").Value = -4.466878511274134
$ws3.Range("I3").Value = -2.729039407369889
$ws3.Range("J3").Value = -10.02112105437275
$ws3.Range("B4").Value = -6.537516359248388
$ws3.Range("C4").Value = 7.878453037160982
$ws3.Range("E4").Value = -6.286917285092488
$ws3.Range("F4").Value = -6.040595615681973
$ws3.Range("G4").Value = -3.679514517577699
$ws3.Range("H4").Value = -3.59713274322795
$ws3.Range("I4").Value = -2.24696183624066
$ws3.Range("J4").Value = -9.730253402761253
$ws3.Range("B5").Value = -5.022538861161706
$ws3.Range("C5").Value = 6.824309757728702
$ws3.Range("D5").Value = 6.286917285092488
$ws3.Range("F5").Value = -4.484345327674603
$ws3.Range("G5").Value = 0.5460780222551406
$ws3.Range("H5").Value = 0.2343698639606305
$ws3.Range("I5").Value = 0.0005651232292951768
$ws3.Range("J5").Value = -6.088110079044103
$ws3.Range("B6").Value = -7.256117192520128
$ws3.Range("C6").Value = 6.331545514136495
$ws3.Range("D6").Value = 6.040595615681973
$ws3.Range("E6").Value = 4.484345327674603
$ws3.Range("G6").Value = 6.597167248742768
$ws3.Range("H6").Value = 6.782378503854368
$ws3.Range("I6").Value = 9.461834552365913
$ws3.Range("J6").Value = 2.585513787940318
$ws3.Range("B7").Value = -7.18893327356635
$ws3.Range("C7").Value = 4.61693374871366
$ws3.Range("D7").Value = 3.679514517577699
$ws3.Range("E7").Value = -0.5460780222551406
$ws3.Range("F7").Value = -6.597167248742768
$ws3.Range("H7").Value = -1.617546462906271
$ws3.Range("I7").Value = -0.4492363056706506
$ws3.Range("J7").Value = -3.860071073428863
$ws3.Range("B8").Value = -7.430047475841675
$ws3.Range("C8").Value = 4.466878511274134
$ws3.Range("D8").Value = 3.59713274322795
$ws3.Range("E8").Value = -0.2343698639606305
$ws3.Range("F8").Value = -6.782378503854368
$ws3.Range("G8").Value = 1.617546462906271
$ws3.Range("I8").Value = -0.2295425250708586
$ws3.Range("J8").Value = -3.395309374597946
$ws3.Range("B9").Value = -10.86844846190304
$ws3.Range("C9").Value = 2.729039407369889
$ws3.Range("D9").Value = 2.24696183624066
$ws3.Range("E9").Value = -0.0005651232292951768
$ws3.Range("F9").Value = -9.461834552365913
$ws3.Range("G9").Value = 0.4492363056706506
$ws3.Range("H9").Value = 0.2295425250708586
$ws3.Range("J9").Value = -2.183738900754142
$ws3.Range("B10").Value = -3.153776194383034
$ws3.Range("C10").Value = 10.02112105437275
$ws3.Range("D10").Value = 9.730253402761253
$ws3.Range("E10").Value = 6.088110079044103
$ws3.Range("F10").Value = -2.585513787940318
$ws3.Range("G10").Value = 3.860071073428863
$ws3.Range("H10").Value = 3.395309374597946
$ws3.Range("I10").Value = 2.183738900754142

# Resumen: 29 cell updates
$ws4.Range("B3").Value = 4
$ws4.Range("D3").Value = 3
$ws4.Range("E3").Value = 50
$ws4.Range("A4").Value = "LSPMW"
$ws4.Range("F4").Value = 2.698648693965452
$ws4.Range("A5").Value = "DeepAR"
$ws4.Range("B5").Value = 2
$ws4.Range("C5").Value = 0
$ws4.Range("D5").Value = 6
$ws4.Range("E5").Value = 25
$ws4.Range("F5").Value = 2.698185999935002
$ws4.Range("B6").Value = 2
$ws4.Range("C6").Value = 1
$ws4.Range("D6").Value = 5
$ws4.Range("E6").Value = 25
$ws4.Range("A7").Value = "AV-MCPS"
$ws4.Range("C7").Value = 1
$ws4.Range("D7").Value = 5
$ws4.Range("F7").Value = 2.579992484281206
$ws4.Range("A9").Value = "Block Bootstrapping"
$ws4.Range("B9").Value = 0
$ws4.Range("C9").Value = 7
$ws4.Range("D9").Value = 1
$ws4.Range("E9").Value = 0
$ws4.Range("F9").Value = 7.880117721526667
$ws4.Range("A10").Value = "EnCQR-LSTM"
$ws4.Range("C10").Value = 3
$ws4.Range("D10").Value = 5
$ws4.Range("F10").Value = 4.570246739591618
